# Rename sheets: strip the "HPOP" prefix from each sheet's name.
#   HPOPdata            -> data
#   HPOPChart           -> Chart
#   HPOPTime Series      -> Time Series
#   HPOPInter           -> Inter
#   HPOPIndicator List  -> Indicator List
#
# Excel automatically repoints ordinary cell formulas (and the Print_Area
# is fixed up separately below since Excel does not auto-update it), chart
# series text.

$wb = $excel.ActiveWorkbook

$dataWs = $wb.Worksheets.Item("HPOPdata")
$chartWs = $wb.Worksheets.Item("HPOPChart")
$timeWs = $wb.Worksheets.Item("HPOPTime Series")
$interWs = $wb.Worksheets.Item("HPOPInter")
$indWs = $wb.Worksheets.Item("HPOPIndicator List")

$dataWs.Name = "data"
$chartWs.Name = "Chart"
$timeWs.Name = "Time Series"
$interWs.Name = "Inter"
$indWs.Name = "Indicator List"

# The named range Print_Area (local to the "data" sheet) keeps pointing at
# the old sheet name after a plain rename, so repoint it explicitly.
$printArea = $wb.Names.Item("data!Print_Area")
$printArea.RefersTo = "=data!`$A`$1:`$Q`$25"

# Repoint the scatter-chart series "name" references (the rest of each
# series' ranges follow the sheet automatically through Excel's normal
# formula-repointing on sheet rename).
$chartObj = $chartWs.ChartObjects(1)
$chart = $chartObj.Chart
$chart.SeriesCollection(1).Name = "=Inter!`$E`$2"
$chart.SeriesCollection(2).Name = "=Inter!`$F`$2"
$chart.SeriesCollection(3).Name = "=Inter!`$I`$2"

# A handful of formulas on the "Inter" sheet reference a deleted column on
# "data" (#REF!). Excel preserves the dead sheet-qualifier text on rename;
# restore that here too.
$interWs.Range("C19").Formula = '=IF(data!#REF!<>"",data!#REF!,#N/A)'
$interWs.Range("D19").Formula = '=IF(data!#REF!<>"",data!#REF!,#N/A)'
$interWs.Range("I19").Formula = '=IF(data!#REF!<>"",data!#REF!,#N/A)'
$interWs.Range("J19").Formula = '=data!#REF!'

# Move the selection/active-tab: the previously-selected "data" sheet's
# cursor moves to B20 and is no longer the active tab; "Indicator List"
# becomes the active tab (cursor stays at G8).
$dataWs.Range("B20").Select() | Out-Null
$indWs.Activate()
$indWs.Range("G8").Select() | Out-Null
